# Updates the "two-digit number divided by one-digit number" practice
# sheet: swaps a batch of division problems for a newly generated set.
#
# Most cells are simple 1-for-1 text substitutions, handled with
# Find/Replace scoped to the whole document. One row (row 9 of the
# table) changed shape in the source edit -- a cell was inserted and two
# trailing cells were collapsed into one -- but the table still has the
# same 5-column grid before and after, so we reproduce that row by
# writing the five target values straight into Table(1).Cell(9, 1..5)
# rather than trying to insert/delete cells.

$d = $word.ActiveDocument

function Replace-Problem($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Simple swaps (order matters only for the one pair noted below).
Replace-Problem "44÷7=" "28÷2="
Replace-Problem "24÷8=" "77÷6="
Replace-Problem "88÷3=" "14÷5="
Replace-Problem "12÷6=" "83÷4="
Replace-Problem "56÷4=" "72÷6="

Replace-Problem "26÷5=" "35÷6="
Replace-Problem "56÷9=" "12÷3="
Replace-Problem "30÷3=" "26÷4="
Replace-Problem "96÷3=" "11÷9="
Replace-Problem "70÷3=" "81÷5="

Replace-Problem "58÷6=" "28÷4="
# "73÷5=" must be replaced BEFORE "61÷6=" is turned into "73÷5=" below,
# otherwise the freshly-written text would get matched and replaced again.
Replace-Problem "73÷5=" "11÷3="
Replace-Problem "10÷3=" "85÷8="
Replace-Problem "69÷7=" "68÷9="
Replace-Problem "90÷9=" "24÷5="

Replace-Problem "62÷3=" "65÷7="
Replace-Problem "63÷2=" "71÷6="
Replace-Problem "61÷6=" "73÷5="
Replace-Problem "24÷9=" "66÷3="
Replace-Problem "40÷3=" "44÷3="

# Row 9: "74÷6=, 73÷2=, 75÷3=, 12÷8=, 35÷7=" becomes
#        "43÷3=, 74÷6=, 89÷8=, 57÷7=, 62÷9="
$t = $d.Tables(1)
$t.Cell(9, 1).Range.Text = "43÷3="
$t.Cell(9, 2).Range.Text = "74÷6="
$t.Cell(9, 3).Range.Text = "89÷8="
$t.Cell(9, 4).Range.Text = "57÷7="
$t.Cell(9, 5).Range.Text = "62÷9="
